$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Đơn sale chính" - populate header row, one data row, one total row
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$headers = @(
    "Tiền tố","Mã dịch vụ","Ngày thực hiện","Cơ sở","Khách hàng","Nguồn khách",
    "Nhóm dịch vụ","Tên dịch vụ","Sale chính","Đơn giá gốc","Sale phụ","Upsale",
    "Đơn giá","Thanh toán lần đầu","Trả sau","Đã thanh toán","Dư nợ","Bác sĩ 1",
    "Bác sĩ 2","Phụ phẫu 1","Phụ phẫu 2","Công phụ phẫu 1","Công phụ phẫu 2",
    "Tỉ lệ chiết khấu sale chính","Tỉ lệ chiết khấu sale phụ","Chiết khấu sale chính",
    "Chiết khấu sale phụ"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws1.Cells.Item(2, 1).Value = "HD-LUXURY"
$ws1.Cells.Item(2, 2).Value = 543
$ws1.Cells.Item(2, 3).Value = "'07-13-2024"
$ws1.Cells.Item(2, 4).Value = "LONG XUYÊN"
$ws1.Cells.Item(2, 5).Value = "Kiều tiên"
$ws1.Cells.Item(2, 6).Value = "Cá nhân"
$ws1.Cells.Item(2, 7).Value = "Tiêm"
$ws1.Cells.Item(2, 8).Value = "Tiêm môi"
$ws1.Cells.Item(2, 9).Value = "Lê Hoàng Thanh"
$ws1.Cells.Item(2, 10).Value = 1800000
$ws1.Cells.Item(2, 13).Value = 1800000
$ws1.Cells.Item(2, 14).Value = 1800000
$ws1.Cells.Item(2, 15).Value = 0
$ws1.Cells.Item(2, 16).Value = 1800000
$ws1.Cells.Item(2, 17).Value = 0
$ws1.Cells.Item(2, 18).Value = "Đặng Ngọc Mai"
$ws1.Cells.Item(2, 22).Value = 0
$ws1.Cells.Item(2, 23).Value = 0
$ws1.Cells.Item(2, 24).Value = 0.1
$ws1.Cells.Item(2, 25).Value = 0
$ws1.Cells.Item(2, 26).Value = 180000
$ws1.Cells.Item(2, 27).Value = 0

$ws1.Cells.Item(3, 1).Value = "Tổng"
$ws1.Cells.Item(3, 2).Value = 1
$ws1.Cells.Item(3, 10).Value = 1800000
$ws1.Cells.Item(3, 12).Value = 0
$ws1.Cells.Item(3, 13).Value = 1800000
$ws1.Cells.Item(3, 14).Value = 1800000
$ws1.Cells.Item(3, 15).Value = 0
$ws1.Cells.Item(3, 16).Value = 1800000
$ws1.Cells.Item(3, 17).Value = 0
$ws1.Cells.Item(3, 22).Value = 0
$ws1.Cells.Item(3, 23).Value = 0
$ws1.Cells.Item(3, 24).Value = 0.1
$ws1.Cells.Item(3, 25).Value = 0
$ws1.Cells.Item(3, 26).Value = 180000
$ws1.Cells.Item(3, 27).Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Lương" - update existing values, insert "Ứng lương" rows and
# append the new "Tổng lương" summary rows
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2, 2).Value = 17
$ws2.Cells.Item(3, 2).Value = 595000

# B4 ("Lương cơ bản tại CẦN THƠ") is a genuinely blank cell (no cached
# value), not zero - make sure it stays that way.
$ws2.Range("B4").ClearContents()

# Insert "Ứng lương tại CẦN THƠ" row right after the CẦN THƠ block (before the
# old row 11 "Lương cơ bản tại LONG XUYÊN"), shifting everything below down by 1.
$ws2.Rows.Item(11).Insert()
$ws2.Cells.Item(11, 1).Value = "Ứng lương tại CẦN THƠ"
$ws2.Cells.Item(11, 2).Value = 0

# Update LONG XUYÊN computed values (now at rows 12-13 after the insert above).
$ws2.Cells.Item(12, 2).Value = 2428571.428571429
$ws2.Cells.Item(13, 2).Value = 180000

# Insert "Ứng lương tại LONG XUYÊN" row right after the LONG XUYÊN block
# (before the row currently holding "Lương cơ bản tại SÓC TRĂNG"), shifting
# everything below down by 1 again.
$ws2.Rows.Item(19).Insert()
$ws2.Cells.Item(19, 1).Value = "Ứng lương tại LONG XUYÊN"
$ws2.Cells.Item(19, 2).Value = 0

# B20 ("Lương cơ bản tại SÓC TRĂNG", shifted down from the original row 18)
# is also genuinely blank - keep it that way.
$ws2.Range("B20").ClearContents()

# Append the trailing rows: "Ứng lương tại SÓC TRĂNG" and the four "Tổng
# lương" totals rows (rows 27-31).
$ws2.Cells.Item(27, 1).Value = "Ứng lương tại SÓC TRĂNG"
$ws2.Cells.Item(27, 2).Value = 0

$ws2.Cells.Item(28, 1).Value = "Tổng lương tại CẦN THƠ"
$ws2.Cells.Item(28, 2).Value = 0

$ws2.Cells.Item(29, 1).Value = "Tổng lương tại LONG XUYÊN"
$ws2.Cells.Item(29, 2).Value = 7666071.428571429

$ws2.Cells.Item(30, 1).Value = "Tổng lương tại SÓC TRĂNG"
$ws2.Cells.Item(30, 2).Value = 0

$ws2.Cells.Item(31, 1).Value = "Tổng lương"
$ws2.Cells.Item(31, 2).Value = 7666071.428571429
